$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $range = $ws.Range($cellAddr)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue "D2" '43.986.56'
Set-TextValue "E2" '  -0.99%  '
Set-TextValue "D3" '2.223.97'
Set-TextValue "E3" '  -0.71%  '
Set-TextValue "D4" '1.00'
Set-TextValue "E4" '  -1.62%  '
Set-TextValue "D5" '298.88'
Set-TextValue "E5" '  -2.65%  '
Set-TextValue "D6" '90.38'
Set-TextValue "E6" '  -4.08%  '
Set-TextValue "E7" '  -2.36%  '
Set-TextValue "D8" '1.00'
Set-TextValue "E8" '  -0.53%  '
Set-TextValue "E9" '  -5.68%  '
Set-TextValue "E10" '  -4.55%  '
Set-TextValue "E11" '  -2.96%  '
Set-TextValue "D12" '6.95'
Set-TextValue "E13" '  -0.68%  '
Set-TextValue "D14" '2.563.10'
Set-TextValue "E14" '  -0.67%  '
Set-TextValue "D15" '2.224.86'
Set-TextValue "E15" '  -1.99%  '
Set-TextValue "D16" '13.40'
Set-TextValue "E16" '  -0.97%  '
Set-TextValue "E17" '  -6.56%  '
Set-TextValue "D18" '43.831.98'
Set-TextValue "E18" '  -0.60%  '
Set-TextValue "E19" '  -0.77%  '
Set-TextValue "D20" '0.0₃0904'
Set-TextValue "E20" '  -4.97%  '
Set-TextValue "D21" '5.96'
Set-TextValue "E21" '  -6.09%  '
Set-TextValue "D22" '64.56'
Set-TextValue "E22" '  -1.40%  '
Set-TextValue "D23" '236.17'
Set-TextValue "E23" '  -0.55%  '
Set-TextValue "D24" '2.81'
Set-TextValue "E24" '  -4.63%  '
Set-TextValue "D25" '0.999'
Set-TextValue "E25" '  -0.35%  '
Set-TextValue "E26" '  -5.77%  '
Set-TextValue "D27" '2.27'
Set-TextValue "E27" '  +2.01%  '
Set-TextValue "D28" '38.91'
Set-TextValue "E28" '  +2.18%  '
Set-TextValue "E29" '  -4.06%  '
Set-TextValue "D30" '151.80'
Set-TextValue "E30" '  -0.77%  '
Set-TextValue "D31" '19.17'
Set-TextValue "E31" '  -3.70%  '
Set-TextValue "E32" '  -9.00%  '
Set-TextValue "D33" '0.0762'
Set-TextValue "E33" '  -3.90%  '
Set-TextValue "E34" '  -5.94%  '
Set-TextValue "E35" '  -1.61%  '
Set-TextValue "D36" '2.84'
Set-TextValue "E36" '  -7.88%  '
Set-TextValue "D37" '0.103'
Set-TextValue "E37" '  -6.85%  '
Set-TextValue "D38" '1.68'
Set-TextValue "E38" '  -6.38%  '
Set-TextValue "D39" '0.0299'
Set-TextValue "E39" '  +0.62%  '
Set-TextValue "D40" '3.61'
Set-TextValue "E40" '  -3.54%  '
Set-TextValue "E41" '  -6.32%  '
Set-TextValue "D42" '13.31'
Set-TextValue "E42" '  -10.81%  '
Set-TextValue "D43" '1.00'
Set-TextValue "E43" '  -0.87%  '
Set-TextValue "D44" '1.801.64'
Set-TextValue "E44" '  +0.64%  '
Set-TextValue "D45" '1.78'
Set-TextValue "E45" '  +13.00%  '
Set-TextValue "E46" '  -3.67%  '
Set-TextValue "D47" '67.68'
Set-TextValue "E47" '  -3.26%  '
Set-TextValue "D48" '94.40'
Set-TextValue "E48" '  -4.17%  '
Set-TextValue "D49" '7.83'
Set-TextValue "E49" '  -2.92%  '
Set-TextValue "D50" '72.87'
Set-TextValue "E50" '  -7.12%  '
Set-TextValue "D51" '4.58'
Set-TextValue "E51" '  -6.11%  '
